$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E15").Value = "see    lcss board bill of materials.xlsx   file"
$ws.Range("F16").Value = "see    lcss programming cable bill of materials.xlsx   file"

$ws.Range("G18").Select()
